# Update countries & provincias Spain
# Applies updated COVID-19 stats and re-orders "Afganistan" into its new
# position (between Azerbaiyan and Camerun) with fresh data, pushing
# Camerun / Bosnia y Herzegovina / Nueva Zelanda down by one row (they
# keep their previous values) and shedding the old "Afganistan" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 11 - Iran
$ws.Range("B11").Value = 90481
$ws.Range("C11").Value = 1153
$ws.Range("D11").Value = 69657
$ws.Range("E11").Value = 15114
$ws.Range("F11").Value = 3079
$ws.Range("G11").Value = 60
$ws.Range("H11").Value = 5710

# Row 52 - Finlandia
$ws.Range("E52").Value = 1886
$ws.Range("F52").Value = 62
$ws.Range("G52").Value = 4
$ws.Range("H52").Value = 190

# Row 72 - Estonia
$ws.Range("B72").Value = 1643
$ws.Range("C72").Value = 8
$ws.Range("D72").Value = 233
$ws.Range("G72").Value = 3
$ws.Range("H72").Value = 49

# Row 74 - becomes Afganistan (new data), was Camerun
$ws.Range("A74").Value = "Afganistan"
$ws.Range("B74").Value = 1531
$ws.Range("C74").Value = 68
$ws.Range("D74").Value = 207
$ws.Range("E74").Value = 1274
$ws.Range("F74").Value = 7
$ws.Range("G74").Value = 3
$ws.Range("H74").Value = 50

# Row 75 - becomes Camerun (old Camerun data), was Bosnia y Herzegovina
$ws.Range("A75").Value = "Camerun"
$ws.Range("B75").Value = 1518
$ws.Range("C75").Value = 0
$ws.Range("D75").Value = 697
$ws.Range("E75").Value = 768
$ws.Range("F75").Value = 28
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 53

# Row 76 - becomes Bosnia y Herzegovina (old Bosnia data), was Nueva Zelanda
$ws.Range("A76").Value = "Bosnia y Herzegovina"
$ws.Range("B76").Value = 1516
$ws.Range("C76").Value = 30
$ws.Range("D76").Value = 624
$ws.Range("E76").Value = 833
$ws.Range("F76").Value = 4
$ws.Range("G76").Value = 2
$ws.Range("H76").Value = 59

# Row 77 - becomes Nueva Zelanda (old Nueva Zelanda data), was Afganistan
$ws.Range("A77").Value = "Nueva Zelanda"
$ws.Range("B77").Value = 1470
$ws.Range("C77").Value = 9
$ws.Range("D77").Value = 1142
$ws.Range("E77").Value = 310
$ws.Range("F77").Value = 1
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 18

# Row 79 - Eslovenia
$ws.Range("B79").Value = 1396
$ws.Range("C79").Value = 8
$ws.Range("D79").Value = 221
$ws.Range("E79").Value = 1093
$ws.Range("G79").Value = 1
$ws.Range("H79").Value = 82

# Row 105 - San Marino
$ws.Range("B105").Value = 538
$ws.Range("C105").Value = 25
$ws.Range("E105").Value = 433
$ws.Range("F105").Value = 4
$ws.Range("G105").Value = 1
$ws.Range("H105").Value = 41

# Row 108 - Sri Lanka
$ws.Range("B108").Value = 471
$ws.Range("C108").Value = 19
$ws.Range("E108").Value = 344

# Row 167 - Nepal
$ws.Range("D167").Value = 16
$ws.Range("E167").Value = 35

# Row 175 - Malaui
$ws.Range("B175").Value = 34
$ws.Range("C175").Value = 1
$ws.Range("E175").Value = 27
